$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (G=5489)
$ws.Range("H2").Value = 328.0909
$ws.Range("I2").Value = 279.85715
$ws.Range("J2").Value = 412.5
$ws.Range("K2").Value = 279.85715
$ws.Range("L2").Value = 412.5
$ws.Range("M2").Value = -166.85715
$ws.Range("N2").Value = -638.5

# Row 5 (G=5503)
$ws.Range("H5").Value = 182.54546
$ws.Range("I5").Value = 184.22223
$ws.Range("J5").Value = 175
$ws.Range("K5").Value = 184.22223
$ws.Range("L5").Value = 175
$ws.Range("M5").Value = -69.22223
$ws.Range("N5").Value = -405

# Row 8 (G=4565)
$ws.Range("H8").Value = 23
$ws.Range("I8").Value = 23
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 69
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = 70

# Row 11 (G=5533)
$ws.Range("H11").Value = 26.285715
$ws.Range("I11").Value = 26.285715
$ws.Range("K11").Value = 26.285715
$ws.Range("M11").Value = 113.714285

# Row 111 (G=27768)
$ws.Range("H111").Value = 719.75
$ws.Range("I111").Value = 751.3333
$ws.Range("J111").Value = 625
$ws.Range("K111").Value = 2253.9999
$ws.Range("L111").Value = 1875
$ws.Range("M111").Value = 813.0001000000002
$ws.Range("N111").Value = -8009

# Row 113 (G=27775)
$ws.Range("H113").Value = 2951
$ws.Range("I113").Value = 3250
$ws.Range("J113").Value = 2353
$ws.Range("K113").Value = 3250
$ws.Range("L113").Value = 2353
$ws.Range("M113").Value = 4
$ws.Range("N113").Value = -8861

# Row 141 (G=44161)
$ws.Range("H141").Value = 1215.381
$ws.Range("I141").Value = 1261.15
$ws.Range("J141").Value = 300
$ws.Range("K141").Value = 3783.45
$ws.Range("L141").Value = 900
$ws.Range("M141").Value = 1396.55
$ws.Range("N141").Value = -11260

$ws = $wb.Worksheets.Item("ARM")
# Row 97 (G=19941)
$ws.Range("H97").Value = 1536.875
$ws.Range("I97").Value = 978.8333
$ws.Range("J97").Value = 3211
$ws.Range("K97").Value = 978.8333
$ws.Range("L97").Value = 3211
$ws.Range("M97").Value = -482.8333
$ws.Range("N97").Value = -4203

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (G=5361)
$ws.Range("H7").Value = 50.95238
$ws.Range("I7").Value = 32
$ws.Range("J7").Value = 81.75
$ws.Range("K7").Value = 32
$ws.Range("L7").Value = 81.75
$ws.Range("M7").Value = 81
$ws.Range("N7").Value = -307.75

# Row 23 (G=2703)
$ws.Range("H23").Value = 8800
$ws.Range("I23").Value = 6733.3335
$ws.Range("J23").Value = 15000
$ws.Range("K23").Value = 6733.3335
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = -6493.3335
$ws.Range("N23").Value = -15480

# Row 27 (G=2703)
$ws.Range("H27").Value = 8800
$ws.Range("I27").Value = 6733.3335
$ws.Range("J27").Value = 15000
$ws.Range("K27").Value = 6733.3335
$ws.Range("L27").Value = 15000
$ws.Range("M27").Value = -6541.3335
$ws.Range("N27").Value = -15384

# Row 31 (G=44023)
$ws.Range("H31").Value = 5403.574
$ws.Range("I31").Value = 3974.8
$ws.Range("K31").Value = 3974.8
$ws.Range("M31").Value = -3679.8

# Row 34 (G=44023)
$ws.Range("H34").Value = 5403.574
$ws.Range("I34").Value = 3974.8
$ws.Range("K34").Value = 3974.8
$ws.Range("M34").Value = -3772.8

# Row 95 (G=18192)
$ws.Range("H95").Value = 17325
$ws.Range("J95").Value = 17325
$ws.Range("L95").Value = 17325
$ws.Range("N95").Value = -22817

# Row 132 (G=44019)
$ws.Range("H132").Value = 2079.5
$ws.Range("I132").Value = 1809.9166
$ws.Range("J132").Value = 3050
$ws.Range("K132").Value = 5429.7498
$ws.Range("L132").Value = 9150
$ws.Range("M132").Value = -2899.7498
$ws.Range("N132").Value = -14210

# Row 134 (G=44020)
$ws.Range("H134").Value = 2665
$ws.Range("I134").Value = 1595.091
$ws.Range("K134").Value = 4785.272999999999
$ws.Range("M134").Value = -2250.272999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 8 (G=16734)
$ws.Range("H8").Value = 408.33334
$ws.Range("I8").Value = 408.33334
$ws.Range("K8").Value = 1225.00002
$ws.Range("M8").Value = -1086.00002

# Row 12 (G=4854)
$ws.Range("H12").Value = 191.64285
$ws.Range("J12").Value = 203.38461
$ws.Range("L12").Value = 610.15383
$ws.Range("N12").Value = -956.15383

# Row 68 (G=12895)
$ws.Range("H68").Value = 734.3333
$ws.Range("J68").Value = 751.5
$ws.Range("L68").Value = 2254.5
$ws.Range("N68").Value = -3876.5

# Row 71 (G=12895)
$ws.Range("H71").Value = 734.3333
$ws.Range("J71").Value = 751.5
$ws.Range("L71").Value = 6763.5
$ws.Range("N71").Value = -14875.5

# Row 103 (G=19839)
$ws.Range("H103").Value = 1578.6428
$ws.Range("J103").Value = 1686.3077
$ws.Range("L103").Value = 5058.9231
$ws.Range("N103").Value = -6816.9231

# Row 116 (G=27866)
$ws.Range("H116").Value = 1950
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 1425
$ws.Range("K116").Value = 9000
$ws.Range("L116").Value = 4275
$ws.Range("M116").Value = -5558
$ws.Range("N116").Value = -11159

# Row 131 (G=36060)
$ws.Range("H131").Value = 1969.0435
$ws.Range("I131").Value = 1425.5
$ws.Range("K131").Value = 4276.5
$ws.Range("M131").Value = 763.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (G=12521)
$ws.Range("H80").Value = 1950
$ws.Range("J80").Value = 1900
$ws.Range("L80").Value = 1900
$ws.Range("N80").Value = -3896

# Row 83 (G=12521)
$ws.Range("H83").Value = 1950
$ws.Range("J83").Value = 1900
$ws.Range("L83").Value = 9500
$ws.Range("N83").Value = -19484

# Row 132 (G=44008)
$ws.Range("H132").Value = 82414.766
$ws.Range("I132").Value = 116466.336
$ws.Range("J132").Value = 5798.75
$ws.Range("K132").Value = 349399.008
$ws.Range("L132").Value = 17396.25
$ws.Range("M132").Value = -346869.008
$ws.Range("N132").Value = -22456.25

$ws = $wb.Worksheets.Item("LTW")
# Row 118 (G=26146)
$ws.Range("H118").Value = 39998
$ws.Range("J118").Value = 39998
$ws.Range("L118").Value = 39998
$ws.Range("N118").Value = -43312

# Row 132 (G=44058)
$ws.Range("H132").Value = 4813.294
$ws.Range("I132").Value = 3485.0833
$ws.Range("J132").Value = 8001
$ws.Range("K132").Value = 10455.2499
$ws.Range("L132").Value = 24003
$ws.Range("M132").Value = -7925.249899999999
$ws.Range("N132").Value = -29063

# Row 136 (G=44060)
$ws.Range("H136").Value = 3957.25
$ws.Range("I136").Value = 2495.6667
$ws.Range("K136").Value = 7487.000100000001
$ws.Range("M136").Value = -4937.000100000001

$ws = $wb.Worksheets.Item("WVR")
# Row 41 (G=21725)
$ws.Range("H41").Value = 40593
$ws.Range("I41").Value = 35997.5
$ws.Range("J41").Value = 45188.5
$ws.Range("K41").Value = 35997.5
$ws.Range("L41").Value = 45188.5
$ws.Range("M41").Value = -35607.5
$ws.Range("N41").Value = -45968.5

# Row 107 (G=27746)
$ws.Range("H107").Value = 875.6667
$ws.Range("I107").Value = 875.6667
$ws.Range("K107").Value = 2627.0001
$ws.Range("M107").Value = -707.0001000000002

# Row 132 (G=44029)
$ws.Range("H132").Value = 3374.375
$ws.Range("I132").Value = 3332.5
$ws.Range("K132").Value = 9997.5
$ws.Range("M132").Value = -7467.5

# Row 136 (G=44031)
$ws.Range("H136").Value = 2187.818
$ws.Range("I136").Value = 1547.5428
$ws.Range("K136").Value = 4642.6284
$ws.Range("M136").Value = -2092.6284
